$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.452.63'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.19'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.06'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4245'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3592'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07214'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8593'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.56'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.780.24'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -5.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.371'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.475'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06931'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.30'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008908'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9999'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.37'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.311.81'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.130'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.36%  '
$ws.Range("E23").Value = '  +2.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.048.95'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.992'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.71'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.63'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.137'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.12'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.785'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -8.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08907'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7464'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.539'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.949'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.122'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.003'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.082'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05252'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01922'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.787'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5067'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1656'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.353'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.376'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.40%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.42'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.35'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06461'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4685'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.002'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.616'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.67'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.64%  '
